# Add a new "Comments" column (column E) with its header to the four
# history sheets (Withdraw History, Deposit History, Transfer History,
# Absolute History). The "Amount" sheet is left untouched.

$wb = $excel.ActiveWorkbook

$historySheetNames = @("Deposit History", "Transfer History", "Absolute History")

foreach ($name in $historySheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "Comments"
    [void]$ws.Range("E1").Select()
}

# "Absolute History" ends up with its selection on E5 rather than E1.
$wsAbsolute = $wb.Worksheets.Item("Absolute History")
[void]$wsAbsolute.Range("E5").Select()

# "Withdraw History" becomes (and stays) the active/selected sheet.
$wsWithdraw = $wb.Worksheets.Item("Withdraw History")
$wsWithdraw.Range("E1").Value = "Comments"
[void]$wsWithdraw.Range("E1").Select()
